$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("logs")

# Row 26: B value changes from the old timestamp to "Pipe_SCTv2_corrected_13-06"
$ws.Range("B26").Value = "Pipe_SCTv2_corrected_13-06"

# Row 29: add G29 = "pseudotime"
$ws.Range("G29").Value = "pseudotime"

# New row 30
$ws.Range("A30").Value = "results"
$ws.Range("B30").Value = "2022-06-14 15-16-04"
$ws.Range("C30").Value = "DEG"
$ws.Range("D30").Value = "SCTv2 corrected BL_N + BL_C new selection"
$ws.Range("F30").Value = "rerun SCTv2 corrected pipeline"
$ws.Range("G30").Value = "pseudotime"

# Update view: scroll position and selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("B30").Select()
